$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update A19 from "led" to "relay"
$ws.Range("A19").Value = "relay"

# Add quantity and price for the new relay row
$ws.Range("B19").Value = 1
$ws.Range("C19").Value = 2.98

# Add the link text for the relay (plain text, not auto-converted to hyperlink)
$ws.Range("D19").Value = "https://www.digikey.com/en/products/detail/sanyou-relay/SRD-S-112DM-F-11/14548486"

# Update the selection to the full row 19 (A19:XFD19), matching the author's last selection
$ws.Range("A19:XFD19").Select()
